# Auto-generated edit script applying numeric corrections to the Kujata Profits workbook.
# Each leve row's price/profit columns (H-N) are updated to the new computed values.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 997.25
$ws.Range("I19").Value = 989
$ws.Range("J19").Value = 998.4286
$ws.Range("K19").Value = 989
$ws.Range("L19").Value = 998.4286
$ws.Range("M19").Value = -814
$ws.Range("N19").Value = -1348.4286
$ws.Range("H40").Value = 2794.6
$ws.Range("I40").Value = 2860.6667
$ws.Range("J40").Value = 2200
$ws.Range("K40").Value = 2860.6667
$ws.Range("L40").Value = 2200
$ws.Range("M40").Value = -2685.6667
$ws.Range("N40").Value = -2550
$ws.Range("H88").Value = 476913.5
$ws.Range("I88").Value = 1540.4286
$ws.Range("J88").Value = 652050.9399999999
$ws.Range("K88").Value = 1540.4286
$ws.Range("L88").Value = 652050.9399999999
$ws.Range("M88").Value = -1134.4286
$ws.Range("N88").Value = -652862.9399999999
$ws.Range("H91").Value = 476913.5
$ws.Range("I91").Value = 1540.4286
$ws.Range("J91").Value = 652050.9399999999
$ws.Range("K91").Value = 1540.4286
$ws.Range("L91").Value = 652050.9399999999
$ws.Range("M91").Value = -136.4286
$ws.Range("N91").Value = -654858.9399999999
$ws.Range("H103").Value = 2760
$ws.Range("I103").Value = 780
$ws.Range("J103").Value = 3750
$ws.Range("K103").Value = 2340
$ws.Range("L103").Value = 11250
$ws.Range("M103").Value = -1754
$ws.Range("N103").Value = -12422
$ws.Range("H113").Value = 2339.5
$ws.Range("I113").Value = 1785.6666
$ws.Range("J113").Value = 4001
$ws.Range("K113").Value = 1785.6666
$ws.Range("L113").Value = 4001
$ws.Range("M113").Value = 1468.3334
$ws.Range("N113").Value = -10509
$ws.Range("H127").Value = 1152.2632
$ws.Range("I127").Value = 555.8125
$ws.Range("J127").Value = 4333.3335
$ws.Range("K127").Value = 1667.4375
$ws.Range("L127").Value = 13000.0005
$ws.Range("M127").Value = 3292.5625
$ws.Range("N127").Value = -22920.0005
$ws.Range("H132").Value = 8555339
$ws.Range("I132").Value = 19616920
$ws.Range("J132").Value = 7754.136
$ws.Range("K132").Value = 58850760
$ws.Range("L132").Value = 23262.408
$ws.Range("M132").Value = -58848230
$ws.Range("N132").Value = -28322.408

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1930.0555
$ws.Range("I74").Value = 1610.0625
$ws.Range("K74").Value = 1610.0625
$ws.Range("M74").Value = -736.0625
$ws.Range("H77").Value = 1930.0555
$ws.Range("I77").Value = 1610.0625
$ws.Range("K77").Value = 8050.3125
$ws.Range("M77").Value = -3682.3125
$ws.Range("H132").Value = 4038.4546
$ws.Range("I132").Value = 3487.4285
$ws.Range("J132").Value = 5002.75
$ws.Range("K132").Value = 10462.2855
$ws.Range("L132").Value = 15008.25
$ws.Range("M132").Value = -7932.2855
$ws.Range("N132").Value = -20068.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 981.8333
$ws.Range("I80").Value = 199.5
$ws.Range("K80").Value = 199.5
$ws.Range("M80").Value = 798.5
$ws.Range("H83").Value = 981.8333
$ws.Range("I83").Value = 199.5
$ws.Range("K83").Value = 997.5
$ws.Range("M83").Value = 3994.5
$ws.Range("H134").Value = 6824.7896
$ws.Range("I134").Value = 1104.5
$ws.Range("J134").Value = 37333
$ws.Range("K134").Value = 3313.5
$ws.Range("L134").Value = 111999
$ws.Range("M134").Value = -778.5
$ws.Range("N134").Value = -117069

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1169.65
$ws.Range("I31").Value = 1088.3077
$ws.Range("J31").Value = 1320.7142
$ws.Range("K31").Value = 1088.3077
$ws.Range("L31").Value = 1320.7142
$ws.Range("M31").Value = -793.3077000000001
$ws.Range("N31").Value = -1910.7142
$ws.Range("H34").Value = 1169.65
$ws.Range("I34").Value = 1088.3077
$ws.Range("J34").Value = 1320.7142
$ws.Range("K34").Value = 1088.3077
$ws.Range("L34").Value = 1320.7142
$ws.Range("M34").Value = -886.3077000000001
$ws.Range("N34").Value = -1724.7142
$ws.Range("H107").Value = 838.8421
$ws.Range("I107").Value = 395.23077
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 395.23077
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 1524.76923
$ws.Range("N107").Value = -5640
$ws.Range("H109").Value = 19500.5
$ws.Range("J109").Value = 19500.5
$ws.Range("L109").Value = 19500.5
$ws.Range("N109").Value = -21580.5
$ws.Range("H132").Value = 6571.4546
$ws.Range("I132").Value = 11562.2
$ws.Range("J132").Value = 2412.5
$ws.Range("K132").Value = 34686.60000000001
$ws.Range("L132").Value = 7237.5
$ws.Range("M132").Value = -32156.60000000001
$ws.Range("N132").Value = -12297.5
$ws.Range("H141").Value = 1015173.3
$ws.Range("J141").Value = 1015173.3
$ws.Range("L141").Value = 1015173.3
$ws.Range("N141").Value = -1025533.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 88.63636
$ws.Range("J2").Value = 88.63636
$ws.Range("L2").Value = 531.81816
$ws.Range("N2").Value = -757.81816
$ws.Range("H32").Value = 2111.3333
$ws.Range("I32").Value = 902
$ws.Range("J32").Value = 2262.5
$ws.Range("K32").Value = 2706
$ws.Range("L32").Value = 6787.5
$ws.Range("M32").Value = -2423
$ws.Range("N32").Value = -7353.5
$ws.Range("H131").Value = 16952170
$ws.Range("I131").Value = 100000310
$ws.Range("J131").Value = 3569.898
$ws.Range("K131").Value = 300000930
$ws.Range("L131").Value = 10709.694
$ws.Range("M131").Value = -299995890
$ws.Range("N131").Value = -20789.694
$ws.Range("H136").Value = 1663.6923
$ws.Range("I136").Value = 1252.9
$ws.Range("K136").Value = 3758.7
$ws.Range("M136").Value = 1341.3

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1464.3334
$ws.Range("I102").Value = 1112
$ws.Range("J102").Value = 1534.8
$ws.Range("K102").Value = 1112
$ws.Range("L102").Value = 1534.8
$ws.Range("M102").Value = 510
$ws.Range("N102").Value = -4778.8
$ws.Range("H113").Value = 1015.0417
$ws.Range("I113").Value = 943.2
$ws.Range("J113").Value = 1134.7778
$ws.Range("K113").Value = 943.2
$ws.Range("L113").Value = 1134.7778
$ws.Range("M113").Value = 1226.8
$ws.Range("N113").Value = -5474.7778
$ws.Range("H114").Value = 29500
$ws.Range("J114").Value = 29500
$ws.Range("L114").Value = 29500
$ws.Range("N114").Value = -38178
$ws.Range("H122").Value = 1630.4
$ws.Range("I122").Value = 1300.5
$ws.Range("J122").Value = 2950
$ws.Range("K122").Value = 3901.5
$ws.Range("L122").Value = 8850
$ws.Range("M122").Value = -1451.5
$ws.Range("N122").Value = -13750

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8262.6875
$ws.Range("I7").Value = 1689.8889
$ws.Range("J7").Value = 16713.428
$ws.Range("K7").Value = 1689.8889
$ws.Range("L7").Value = 16713.428
$ws.Range("M7").Value = -1577.8889
$ws.Range("N7").Value = -16937.428
$ws.Range("H40").Value = 2606.4285
$ws.Range("I40").Value = 1811.7084
$ws.Range("J40").Value = 7374.75
$ws.Range("K40").Value = 1811.7084
$ws.Range("L40").Value = 7374.75
$ws.Range("M40").Value = -1675.7084
$ws.Range("N40").Value = -7646.75
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H100").Value = 1750.75
$ws.Range("I100").Value = 1534.3334
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 1534.3334
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -993.3334
$ws.Range("N100").Value = -3482
$ws.Range("H126").Value = 8262.6875
$ws.Range("I126").Value = 1689.8889
$ws.Range("J126").Value = 16713.428
$ws.Range("K126").Value = 5069.6667
$ws.Range("L126").Value = 50140.284
$ws.Range("M126").Value = -2599.6667
$ws.Range("N126").Value = -55080.284

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 52638990
$ws.Range("I126").Value = 62502244
$ws.Range("J126").Value = 35001.332
$ws.Range("K126").Value = 187506732
$ws.Range("L126").Value = 105003.996
$ws.Range("M126").Value = -187504262
$ws.Range("N126").Value = -109943.996
